# Replace the "Perseus <dates for 2018>" observation-period blurb with the
# translated Cygnus dates, collapsing the whole run-soup (and, where present,
# the trailing "Ennen kuin menet ulos..." sentence + hyperlink) down to a
# single plain run per paragraph.
$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "Perseus havainnointijaksot vuonna 2018*") {
        $r = $p.Range
        # Exclude the trailing paragraph mark so we only touch run content.
        $r.End = $r.End - 1
        $r.Delete()
        $r.InsertAfter("havainnointijaksot vuonna Cygnus: 10.-19.")
    }
}
